$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The sheet currently has a header row (1) and one example/template
# data row (2). The edit pushes that example row down to row 3 (so it
# stays as reference/template data) and replaces row 2 with a new,
# real data row ("满帮" / 雨花区万博科技园 / ...).
# ------------------------------------------------------------------

# Step 1: duplicate row 2 (content + formatting) down into row 3 so the
# original example row is preserved there.
$ws.Range("A2:P2").Copy()
$ws.Range("A3:P3").PasteSpecial(-4104)
$excel.CutCopyMode = 0

# Step 2: row-insert/copy can perturb the style of A3 (it should stay
# identical to A2's style). Re-apply just the formatting from A2 onto
# A3 with a formats-only paste so it keeps the same style id.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Step 3: bump the sequence number that used to be in A2 (0) to 1 now
# that the example row lives in row 3.
$ws.Range("A3").Value = 1

# Step 4: overwrite row 2 with the new company's info.
$ws.Range("B2").Value = "满帮"
$ws.Range("C2").Value = "雨花区万博科技园"

$ws.Range("D2").ClearContents()
$ws.Range("D2").Style = "Normal"

$ws.Range("H2").Value = "看部门，不强制， 周五基本不加，还有每月一天奋斗日（年底算工资）， 据说要取消了"
$ws.Range("I2").Value = "全额8%"
$ws.Range("J2").Value = "上下半年绩效"

$ws.Range("K2").ClearContents()
$ws.Range("K2").Style = "Normal"

$ws.Range("L2").Value = "联想"

$ws.Range("M2").ClearContents()
$ws.Range("M2").Style = "Normal"
$ws.Range("N2").ClearContents()
$ws.Range("N2").Style = "Normal"

# Step 5: row 3's trailing columns (O3/P3) need to exist as blank
# cells, same as O2/P2 in row 2 - materialize them.
$ws.Range("O3").Style = "Normal"
$ws.Range("P3").Style = "Normal"

Write-Output "done"
